$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (class label swapped from M -> B, with new metric values)
$ws.Range("A2").Value = "B"
$ws.Range("B2").Value = 0.9067796610169492
$ws.Range("C2").Value = 0.9953488372093023
$ws.Range("D2").Value = 0.9490022172949002
$ws.Range("E2").Value = 215

# Row 3 (class label swapped from B -> M, with new metric values)
$ws.Range("A3").Value = "M"
$ws.Range("B3").Value = 0.9905660377358491
$ws.Range("C3").Value = 0.8267716535433071
$ws.Range("D3").Value = 0.9012875536480687
$ws.Range("E3").Value = 127

# Row 4 (accuracy)
$ws.Range("B4").Value = 0.9327485380116959
$ws.Range("C4").Value = 0.9327485380116959
$ws.Range("D4").Value = 0.9327485380116959
$ws.Range("E4").Value = 0.9327485380116959

# Row 5 (macro avg)
$ws.Range("B5").Value = 0.9486728493763992
$ws.Range("C5").Value = 0.9110602453763047
$ws.Range("D5").Value = 0.9251448854714844
$ws.Range("E5").Value = 342

# Row 6 (weighted avg)
$ws.Range("B6").Value = 0.9378933155295232
$ws.Range("C6").Value = 0.9327485380116959
$ws.Range("D6").Value = 0.9312836141278019
$ws.Range("E6").Value = 342
